$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: keep username/password, change role from "admin" to "editor"
$ws.Range("C2").Value = "editor"

# Row 3: replace john/pass123/editor with maulik/te#ch$1234/admin
$ws.Range("A3").Value = "maulik"
$ws.Range("B3").Value = "te#ch`$1234"
$ws.Range("C3").Value = "admin"

# Row 4: remove the alice/123/viewer record entirely (leaves an empty row, not a shift)
$ws.Range("A4:C4").ClearContents()

# Update the active selection to C2, matching the saved view state
$ws.Range("C2").Select()
